# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated data, as produced at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Row number (in column F) -> new value
$updates = @{
    2  = 8840
    3  = 8223
    8  = 144
    9  = 154
    10 = 205
    12 = 744
    14 = 5307
    16 = 82
    17 = 21
    20 = 146
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
